$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 22 de Octubre de 2020 a las 18:07"

$rowUpdates = @(
  @{ Row = 4; Values = @("Estados Unidos", 8595023, 10204, 5612505, 2755002, 0, 107, 227516) },
  @{ Row = 5; Values = @("India", 7727289, 22131, 6903365, 706913, 0, 358, 117011) },
  @{ Row = 6; Values = @("Brasil", 5303520, 2871, 4756489, 391531, 0, 41, 155500) },
  @{ Row = 17; Values = @("Chile", 497131, 1494, 469765, 13574, 0, 73, 13792) },
  @{ Row = 18; Values = @("Italia", 465726, 16079, 259456, 169302, 0, 136, 36968) },
  @{ Row = 20; Values = @("Alemania", 395442, 4087, 302100, 83324, 0, 19, 10018) },
  @{ Row = 21; Values = @("Banglades", 394827, 1696, 310532, 78548, 0, 24, 5747) },
  @{ Row = 33; Values = @("Canada", 208234, 2280, 175255, 23121, 0, 32, 9858) },
  @{ Row = 41; Values = @("Republica Dominicana", 122873, 475, 100920, 19741, 0, 6, 2212) },
  @{ Row = 46; Values = @("Portugal", 109541, 3270, 64531, 42765, 0, 16, 2245) },
  @{ Row = 47; Values = @("Suecia", 108969, 0, 0, 0, 0, 6, 5930) },
  @{ Row = 72; Values = @("Estado de Palestina", 49134, 506, 42544, 6155, 0, 8, 435) },
  @{ Row = 73; Values = @("Ghana", 47538, 77, 46789, 437, 0, 0, 312) },
  @{ Row = 75; Values = @("Kenia", 47212, 1068, 33050, 13292, 0, 12, 870) },
  @{ Row = 76; Values = @("Jordania", 46441, 2821, 7340, 38620, 0, 38, 481) },
  @{ Row = 77; Values = @("Tunez", 45892, 0, 5032, 40120, 0, 0, 740) },
  @{ Row = 87; Values = @("Grecia", 28216, 882, 9989, 17678, 0, 15, 549) },
  @{ Row = 88; Values = @("Australia", 27466, 22, 25159, 1402, 0, 0, 905) },
  @{ Row = 99; Values = @("Montenegro", 16259, 190, 12093, 3913, 0, 3, 253) },
  @{ Row = 105; Values = @("Luxemburgo", 12333, 662, 8474, 3719, 0, 2, 140) },
  @{ Row = 122; Values = @("Cuba", 6421, 53, 5871, 422, 0, 1, 128) }
)

foreach ($update in $rowUpdates) {
  $rowNum = $update.Row
  $vals = $update.Values
  for ($c = 0; $c -lt $vals.Length; $c++) {
    $ws.Cells.Item($rowNum, $c + 1).Value = $vals[$c]
  }
}
